$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 23-32: update column C (content) and, where needed, column B (type).
# Column A (keyword) and D (userId) are unchanged for these rows.

$ws.Range("B23").Value = "Significant Related Persons"
$ws.Range("C23").Value = "brian_armstrong"

$ws.Range("B24").Value = "Significant Related Persons"
$ws.Range("C24").Value = "CoinDesk"

$ws.Range("B25").Value = "Significant Related Persons"
$ws.Range("C25").Value = "CryptosisAI"

$ws.Range("B26").Value = "Associated Technology"
$ws.Range("C26").Value = "Serenity"

$ws.Range("B27").Value = "Associated Technology"
$ws.Range("C27").Value = "Casper"

$ws.Range("B28").Value = "Associated Technology"
$ws.Range("C28").Value = "VDF"

$ws.Range("B29").Value = "Associated Technology"
$ws.Range("C29").Value = "Beacon chain"

$ws.Range("B30").Value = "Associated Technology"
$ws.Range("C30").Value = "ERC"

$ws.Range("B31").Value = "Significant Related Persons"
$ws.Range("C31").Value = "ForbesCrypto"

$ws.Range("B32").Value = "Influential Event"
$ws.Range("C32").Value = "merge"

$ws.Range("B33").Value = "Associated Technology"
$ws.Range("C33").Value = "blockchain"
$ws.Range("E33").Value = "The core technology of Ethereum"

$ws.Range("B34").Value = "Significant Related Persons"
$ws.Range("C34").Value = "VitalikButerin"
$ws.Range("E34").Value = "Russian-born Canadian programmer and writer who is best known as one of the co-founders of Ethereum"

$ws.Range("B35").Value = "Associated Technology"
$ws.Range("C35").Value = "smart contracts"
$ws.Range("E35").Value = "A smart contract, like any contract, establishes the terms of an agreement. But unlike a traditional contract, a smart contract’s terms are executed as code running on a blockchain like Ethereum. Smart contracts allow developers to build apps that take advantage of blockchain security, reliability, and accessibility while offering sophisticated peer-to-peer functionality — everything from loans and insurance to logistics and gaming."

$ws.Range("B36").Value = "Correlated Concept"
$ws.Range("C36").Value = "ETC"
$ws.Range("E36").Value = "ETC is the native cryptocurrency of Ethereum Classic, a blockchain project that was created in 2016 when Ethereum’s blockchain split into two separate chains following a disagreement among members of its community."

$ws.Range("B37").Value = "Correlated Concept"
$ws.Range("C37").Value = "DeFi"
$ws.Range("E37").Value = "Decentralized finance (DeFi) is an emerging financial technology based on secure distributed ledgers similar to those used by cryptocurrencies. The system removes the control banks and institutions have on money, financial products, and financial services."

# New row 38
$ws.Range("A38").Value = "eth"
$ws.Range("B38").Value = "Correlated Concept"
$ws.Range("C38").Value = "merge"
$ws.Range("D38").Value = "test111"
$ws.Range("E38").Value = "Eventually the current Ethereum Mainnet will `"merge`" with the beacon chain proof-of-stake system.`nThis will mark the end of proof-of-work for Ethereum, and the full transition to proof-of-stake.`nThis is planned to precede the roll out of shard chains.`nWe formerly referred to this as `"the docking.`""
